# LoginData.xlsx update
#  - doLogin!D1 header renamed "expectedTest" -> "expectedTitle"; selection moved to D1
#  - new sheet "logedinToDashboard" (copy of doLogin, same data/format)
#  - new sheet "addTestFromDropDown" (username/password/browserName/testName/discount)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- doLogin: rename header, move selection -------------------------------
$ws1.Range("D1").Value = "expectedTitle"
$ws1.Range("D1").Select() | Out-Null

# --- logedinToDashboard: duplicate of doLogin ------------------------------
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$s2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s2.Name = "logedinToDashboard"

# --- addTestFromDropDown: new sheet with test/discount data ---------------
$s3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$s3.Name = "addTestFromDropDown"

$s3.Range("A1").Value = "username"
$s3.Range("B1").Value = "password"
$s3.Range("C1").Value = "browserName"
$s3.Range("D1").Value = "testName"
$s3.Range("E1").Value = "discount"

$s3.Range("A2").Value = "test@kennect.io"
$s3.Range("B2").Value = "Qwerty@1234"
$s3.Range("C2").Value = "chrome"
$s3.Range("D2").Value = "Beans"
$s3.Range("E2").Value = "none"

$s3.Range("A3").Value = "test@kennect.io"
$s3.Range("B3").Value = "Qwerty@1234"
$s3.Range("C3").Value = "edge"
$s3.Range("D3").Value = "xyz"
$s3.Range("E3").Value = 0.05

$s3.Range("A4").Value = "test@kennect.io"
$s3.Range("B4").Value = "Qwerty@1234"
$s3.Range("C4").Value = "firefox"
$s3.Range("D4").Value = "UR Uric acid"
$s3.Range("E4").Value = 0.1

$s3.Range("A5").Value = "test@kennect.io"
$s3.Range("B5").Value = "Qwerty@1234"
$s3.Range("C5").Value = "chrome"
$s3.Range("D5").Value = "xray"
$s3.Range("E5").Value = 0.05

$s3.Range("E2:E5").NumberFormat = "0%"

$s3.Hyperlinks.Add($s3.Range("A2"), "mailto:test@kennect.io") | Out-Null
$s3.Hyperlinks.Add($s3.Range("B2"), "mailto:Qwerty@1234") | Out-Null
$s3.Hyperlinks.Add($s3.Range("A4"), "mailto:test@kennect.io") | Out-Null

# Hyperlinks.Add() overwrites the cell style; restore the plain Hyperlink
# style (same one doLogin already uses) on the linked cells.
$s3.Range("A2:B5").Style = "Hyperlink"
$s3.Range("E7").Select() | Out-Null

# Keep doLogin as the active/selected tab (copying a sheet activates the copy).
$ws1.Activate() | Out-Null
